$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# Insert a new row before row 14, shifting rows 14-19 down to 15-20
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 (new shared string "Simple_product_1qty"
# must be created before "vegan-silk-scarf" to match the target string order)
# Note: Rows().Insert() already copies formatting from the row above, so V14/Y14
# already carry the correct styles.
$ws.Range("A14").Value = "Simple_product_1qty"

# Update product names referenced in column U for rows 4, 5, 7 and (now) 13-16
$ws.Range("U4").Value = "Hold Me Softly Style Balm"
$ws.Range("U5").Value = "Hold Me Softly Style Balm"
$ws.Range("U7").Value = "Mini 2 fl oz."
$ws.Range("U13").Value = "vegan-silk-scarf"
$ws.Range("U14").Value = "vegan-silk-scarf"
$ws.Range("U15").Value = "vegan-silk-scarf"
$ws.Range("U16").Value = "vegan-silk-scarf"

# The row insertion shifted hyperlink anchor cells (K16->K17, B18->B19, D18->D19)
# but does not auto-update the worksheet Hyperlinks collection, so rebuild it:
# drop every existing hyperlink and recreate them anchored at their new cells.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Lotuswave@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:vnarra@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:avayugundla@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Lotuswave@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K17"), "mailto:vnarra@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B19"), "mailto:avayugundla@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D19"), "mailto:Lotuswave@123") | Out-Null

# Update the selected cell in the sheet view
$ws.Range("C8").Select()
